$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet_name_3"

# Unicode right single quotation mark used in several strings
$q = [char]0x2018

# Row 2 (data index 0): clear B2 and C2 entirely
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

# Row 3 (data index 1)
$ws.Range("B3").Value2 = "canada"
$ws.Range("C3").Value2 = "3205775 Canagian robots "

# Row 4 (data index 2)
$ws.Range("B4").Value2 = "  "
$ws.Range("C4").Value2 = "  "

# Row 5 (data index 3)
$ws.Range("B5").Value2 = "Ohio"
$ws.Range("C5").Value2 = "kang Matartic  "

# Row 6 (data index 4): clear B6 and C6 entirely
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Row 7 (data index 5): clear B7 and C7 entirely
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 8 (data index 6): clear B8 and C8 entirely
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

# Row 9 (data index 7)
$ws.Range("B9").Value2 = "embados"
$ws.Range("C9").Value2 = "Penns Banking  "

# Row 10 (data index 8)
$ws.Range("B10").Value2 = "Chie"
$ws.Range("C10").Value2 = "Animas Transaction Ine. "

# Row 11 (data index 9)
$ws.Range("B11").Value2 = "Ohio"
$ws.Range("C11").Value2 = "2azi4s1 " + $q + "hie ine. "

# Row 12 (data index 10)
$ws.Range("B12").Value2 = "Urastar " + $q + "oppurinity " + $q + "Corp.  Gries columbia"
$ws.Range("C12").Value2 = "Urastar " + $q + "oppurinity " + $q + "Corp.  Gries columbia) "

# Row 13 (data index 11)
$ws.Range("B13").Value2 = "Ohio"
$ws.Range("C13").Value2 = "      2azi4s1 " + $q + "hie ine. "

# Row 14 (data index 12) - new row
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = "Ohio"
$ws.Range("C14").Value2 = "Anico El MINES LIMITED "

# Row 15 (data index 13) - new row
$ws.Range("A13").Copy($ws.Range("A15"))
$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = "Ohio"
$ws.Range("C15").Value2 = "      Anico El MINES LIMITED "
